$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '35.060.60'
$ws.Range('E2').Value = '  +1.63%  '

# Row 3
$ws.Range('D3').Value = '1.858.48'
$ws.Range('E3').Value = '  +3.36%  '

# Row 4
$ws.Range('E4').Value = '  +0.14%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '238.07'
$ws.Range('E5').Value = '  +4.04%  '

# Row 6
$ws.Range('E6').Value = '  +1.85%  '

# Row 7
$ws.Range('E7').Value = '  +0.09%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '42.35'
$ws.Range('E8').Value = '  +9.16%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.328'
$ws.Range('E9').Value = '  +2.96%  '

# Row 10
$ws.Range('E10').Value = '  +3.11%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0991'
$ws.Range('E11').Value = '  +0.39%  '

# Row 12
$ws.Range('D12').Value = '2.128.75'
$ws.Range('E12').Value = '  +3.33%  '

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.860.25'
$ws.Range('E13').Value = '  +3.25%  '

# Row 14
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '11.39'
$ws.Range('E14').Value = '  +2.86%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.678'
$ws.Range('E15').Value = '  +3.29%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '4.69'
$ws.Range('E16').Value = '  +3.55%  '

# Row 17
$ws.Range('D17').Value = '35.051.08'
$ws.Range('E17').Value = '  +1.59%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '70.32'
$ws.Range('E18').Value = '  +2.15%  '

# Row 19
$ws.Range('E19').Value = '  +2.78%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '240.92'
$ws.Range('E20').Value = '  +0.71%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.17'
$ws.Range('E21').Value = '  +3.74%  '

# Row 22
$ws.Range('E22').Value = '  +1.63%  '

# Row 23
$ws.Range('E23').Value = '  -0.08%  '

# Row 24
$ws.Range('E24').Value = '  +2.24%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '171.56'
$ws.Range('E25').Value = '  -0.34%  '

# Row 26
$ws.Range('E26').Value = '  +31.03%  '

# Row 27
$ws.Range('E27').Value = '  +3.16%  '

# Row 28
$ws.Range('E28').Value = '  +3.43%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.124'
$ws.Range('E29').Value = '  +3.29%  '

# Row 30
$ws.Range('B30').Value = 'BinanceUSD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.01'
$ws.Range('E30').Value = '  +0.19%  '

# Row 31
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0558'
$ws.Range('E31').Value = '  +3.55%  '

# Row 32
$ws.Range('E32').Value = '  -0.40%  '

# Row 33
$ws.Range('E33').Value = '  +3.25%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.02'
$ws.Range('E34').Value = '  +13.69%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.62'
$ws.Range('E35').Value = '  +22.87%  '

# Row 36
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.788'
$ws.Range('E36').Value = '  +14.41%  '

# Row 37
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.30'
$ws.Range('E37').Value = '  +5.63%  '

# Row 38
$ws.Range('E38').Value = '  +12.88%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '91.66'
$ws.Range('E39').Value = '  +1.10%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0203'
$ws.Range('E40').Value = '  +6.90%  '

# Row 41
$ws.Range('D41').Value = '1.353.19'
$ws.Range('E41').Value = '  +3.12%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '14.85'
$ws.Range('E42').Value = '  +4.91%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.35'
$ws.Range('E43').Value = '  +6.60%  '

# Row 44
$ws.Range('B44').Value = 'HuobiToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.41'
$ws.Range('E44').Value = '  -0.56%  '

# Row 45
$ws.Range('B45').Value = 'Gas'
$ws.Range('C45').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '12.52'
$ws.Range('E45').Value = '  +54.73%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.75'
$ws.Range('E46').Value = '  +1.93%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0548'
$ws.Range('E47').Value = '  +6.81%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '6.35'
$ws.Range('E48').Value = '  +2.73%  '

# Row 49
$ws.Range('D49').Value = '2.040.30'
$ws.Range('E49').Value = '  +2.97%  '

# Row 50
$ws.Range('E50').Value = '  +3.21%  '

# Row 51
$ws.Range('E51').Value = '  +17.90%  '
